# Applies the crypto price/volume refresh described in the commit diff.
# Column D ("Price") values that parse as plain numbers get an apostrophe
# text-prefix (then Style reset to "Normal") so they stay text cells, like
# the source data (Column D is stored as inline strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.103.75"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3
$ws.Range("D3").Value = "1.681.92"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'215.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6
$ws.Range("D6").Value = "'0.518"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  +2.12%  "

# Row 9
$ws.Range("D9").Value = "'21.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.89%  "

# Row 10
$ws.Range("D10").Value = "'0.0623"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

# Row 11
$ws.Range("E11").Value = "  -0.58%  "

# Row 12
$ws.Range("D12").Value = "1.916.53"
$ws.Range("E12").Value = "  +0.46%  "

# Row 13
$ws.Range("D13").Value = "1.685.37"
$ws.Range("E13").Value = "  +3.11%  "

# Row 14
$ws.Range("D14").Value = "'4.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.61%  "

# Row 15
$ws.Range("E15").Value = "  +2.15%  "

# Row 16
$ws.Range("D16").Value = "'66.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.86%  "

# Row 17
$ws.Range("D17").Value = "27.090.82"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$ws.Range("D18").Value = "'238.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.59%  "

# Row 19
$ws.Range("E19").Value = "  +0.73%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0748"
$ws.Range("E20").Value = "  +2.14%  "

# Row 21
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
$ws.Range("D22").Value = "'4.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.64%  "

# Row 23
$ws.Range("D23").Value = "'9.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.78%  "

# Row 25
$ws.Range("D25").Value = "'146.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "

# Row 26
$ws.Range("E26").Value = "  +0.96%  "

# Row 27
$ws.Range("D27").Value = "'16.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.25%  "

# Row 28
$ws.Range("E28").Value = "  +0.85%  "

# Row 29
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("E30").Value = "  +0.45%  "

# Row 31
$ws.Range("E31").Value = "  +0.26%  "

# Row 32
$ws.Range("D32").Value = "1.560.04"
$ws.Range("E32").Value = "  +6.02%  "

# Row 33
$ws.Range("E33").Value = "  +1.10%  "

# Row 34
$ws.Range("D34").Value = "'3.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.94%  "

# Row 35
$ws.Range("E35").Value = "  +2.40%  "

# Row 36
$ws.Range("D36").Value = "'0.606"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.99%  "

# Row 37
$ws.Range("E37").Value = "  +4.87%  "

# Row 39
$ws.Range("E39").Value = "  +2.23%  "

# Row 40
$ws.Range("E40").Value = "  +0.33%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'68.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.29%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "

# Row 43
$ws.Range("D43").Value = "'5.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "

# Row 44
$ws.Range("D44").Value = "'2.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.81%  "

# Row 45
$ws.Range("D45").Value = "1.825.29"
$ws.Range("E45").Value = "  +0.59%  "

# Row 46
$ws.Range("E46").Value = "  +0.55%  "

# Row 47
$ws.Range("D47").Value = "'90.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("D48").Value = "'1.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "

# Row 49
$ws.Range("E49").Value = "  +0.86%  "

# Row 50
$ws.Range("D50").Value = "'0.105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.66%  "

# Row 51
$ws.Range("D51").Value = "'8.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.35%  "
